$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(55, 'he; boyfriend', '彼|かれ'),
    @(56, 'she; girlfriend', '彼女|かのじょ'),
    @(57, 'they', '彼ら|かれら'),
    @(58, 'boyfriend', '彼氏|かれし'),
    @(59, 'age; era', '時代|じだい'),
    @(60, 'electricity fee', '電気代|でんきだい'),
    @(61, '90''s', '九十年代|きゅうじゅうねんだい'),
    @(62, 'in one''s teens', '十代|じゅうだい'),
    @(63, 'instead', '代わりに|かわりに'),
    @(64, 'international students', '留学生|りゅうがくせい'),
    @(65, 'to study abroad', '留学する|りゅうがくする'),
    @(66, 'absence; not home', '留守|るす'),
    @(67, 'family', '家族|かぞく'),
    @(68, 'race; ethnic group', '民族|みんぞく'),
    @(69, 'aquarium', '水族館|すいぞくかん'),
    @(70, 'member of royalty', '王族|おうぞく'),
    @(71, 'father', '父親|ちちおや'),
    @(72, 'kind', '親切な|しんせつな'),
    @(73, 'best friend', '親友|しんゆう'),
    @(74, 'parents', '両親|りょうしん'),
    @(75, 'intimate', '親しい|したしい'),
    @(76, 'mother', '母親|ははおや'),
    @(77, 'to cut', '切る|きる'),
    @(78, 'ticket', '切符|きっぷ'),
    @(79, 'postage stamp', '切手|きって'),
    @(80, 'precious', '大切な|たいせつな'),
    @(81, 'English language', '英語|えいご'),
    @(82, 'United Kingdom', '英国|えいこく'),
    @(83, 'English conversation', '英会話|えいかいわ'),
    @(84, 'hero', '英雄|えいゆう'),
    @(85, 'shop', '店|みせ'),
    @(86, 'store clerk', '店員|てんいん'),
    @(87, 'stall; kiosk', '売店|ばいてん'),
    @(88, 'book store', '書店|しょてん'),
    @(89, 'store manager', '店長|てんちょう'),
    @(90, 'last year', '去年|きょねん'),
    @(91, 'the past', '過去|かこ'),
    @(92, 'to leave', '去る|さる'),
    @(93, 'to erase', '消去する|しょうきょする'),
    @(94, 'suddenly', '急に|きゅうに'),
    @(95, 'to hurry', '急ぐ|いそぐ'),
    @(96, 'express train', '急行|きゅうこう'),
    @(97, 'super express', '特急|とっきゅう'),
    @(98, 'to ride', '乗る|のる'),
    @(99, 'vehicle', '乗り物|のりもの'),
    @(100, 'riding a car', '乗車|じょうしゃ'),
    @(101, 'horseback riding', '乗馬|じょうば'),
    @(102, 'really', '本当に|ほんとうに'),
    @(103, 'lunch box', 'お弁当|おべんとう'),
    @(104, 'at that time', '当時|とうじ'),
    @(105, 'to hit', '当たる|あたる'),
    @(106, 'music', '音楽|おんがく'),
    @(107, 'pronunciation', '発音|はつおん'),
    @(108, 'sound', '音|おと'),
    @(109, 'real intention', '本音|ほんね'),
    @(110, 'fun', '楽しい|たのしい'),
    @(111, 'musical instrument', '楽器|がっき'),
    @(112, 'easy; comfortable', '楽な|らくな'),
    @(113, 'doctor', '医者|いしゃ'),
    @(114, 'dentist', '歯医者|はいしゃ'),
    @(115, 'medical science', '医学|いがく'),
    @(116, 'clinic', '医院|いいん'),
    @(117, 'scholar', '学者|がくしゃ'),
    @(118, 'reader', '読者|どくしゃ'),
    @(119, 'young people', '若者|わかもの'),
    @(120, 'ninja', '忍者|にんじゃ'),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}